# Auto-generated edit script: updates cryptos list values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.319.15'
$ws.Range('E2').Value = '  +0.65%  '

$ws.Range('D3').Value = '3.343.31'
$ws.Range('E3').Value = '  +0.19%  '

$ws.Range('E4').Value = '  -0.06%  '

$ws.Range('D5').Value = '''556.42'
$ws.Range('E5').Value = '  +0.05%  '

$ws.Range('D6').Value = '''174.52'
$ws.Range('E6').Value = '  -0.19%  '

$ws.Range('D7').Value = '''0.618'
$ws.Range('E7').Value = '  +0.93%  '

$ws.Range('D8').Value = '3.339.06'
$ws.Range('E8').Value = '  +0.23%  '

$ws.Range('E9').Value = '  -0.01%  '

$ws.Range('E10').Value = '  +3.61%  '

$ws.Range('E11').Value = '  +1.69%  '

$ws.Range('D12').Value = '''53.72'
$ws.Range('E12').Value = '  -0.41%  '

$ws.Range('D13').Value = '''0.0000274'
$ws.Range('E13').Value = '  +1.24%  '

$ws.Range('D14').Value = '''9.06'
$ws.Range('E14').Value = '  +0.69%  '

$ws.Range('D15').Value = '3.879.63'
$ws.Range('E15').Value = '  -0.01%  '

$ws.Range('E16').Value = '  +1.15%  '

$ws.Range('D17').Value = '''18.13'
$ws.Range('E17').Value = '  -0.88%  '

$ws.Range('D18').Value = '3.351.63'
$ws.Range('E18').Value = '  +0.12%  '

$ws.Range('D19').Value = '64.725.30'
$ws.Range('E19').Value = '  +1.29%  '

$ws.Range('D20').Value = '''11.71'
$ws.Range('E20').Value = '  -1.06%  '

$ws.Range('D21').Value = '''0.986'
$ws.Range('E21').Value = '  +0.86%  '

$ws.Range('D22').Value = '''452.78'
$ws.Range('E22').Value = '  +4.37%  '

$ws.Range('D23').Value = '''4.94'
$ws.Range('E23').Value = '  +8.45%  '

$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D24').Value = '''14.34'
$ws.Range('E24').Value = '  +8.99%  '

$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').Value = '''4.10'
$ws.Range('E25').Value = '  -0.44%  '

$ws.Range('D26').Value = '''87.09'
$ws.Range('E26').Value = '  +3.36%  '

$ws.Range('D27').Value = '''2.88'
$ws.Range('E27').Value = '  +1.84%  '

$ws.Range('D28').Value = '''10.74'
$ws.Range('E28').Value = '  +1.00%  '

$ws.Range('D29').Value = '''8.70'
$ws.Range('E29').Value = '  -0.48%  '

$ws.Range('D30').Value = '''31.26'
$ws.Range('E30').Value = '  +6.08%  '

$ws.Range('D31').Value = '''6.57'
$ws.Range('E31').Value = '  +0.69%  '

$ws.Range('D32').Value = '''11.46'
$ws.Range('E32').Value = '  +0.03%  '

$ws.Range('D33').Value = '''569.33'
$ws.Range('E33').Value = '  -3.80%  '

$ws.Range('D34').Value = '''60.95'
$ws.Range('E34').Value = '  +3.98%  '

$ws.Range('D35').Value = '''0.107'
$ws.Range('E35').Value = '  +0.15%  '

$ws.Range('E36').Value = '  -0.03%  '

$ws.Range('D37').Value = '''3.64'
$ws.Range('E37').Value = '  +4.00%  '

$ws.Range('D38').Value = '''0.140'
$ws.Range('E38').Value = '  -2.45%  '

$ws.Range('D39').Value = '''35.28'
$ws.Range('E39').Value = '  -0.13%  '

$ws.Range('B40').Value = 'TheGraph'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D40').Value = '''0.368'
$ws.Range('E40').Value = '  +0.82%  '

$ws.Range('B41').Value = 'PEPE'
$ws.Range('C41').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D41').Value = '0.0₃0736'
$ws.Range('E41').Value = '  -1.70%  '

$ws.Range('D42').Value = '3.059.65'
$ws.Range('E42').Value = '  -1.68%  '

$ws.Range('D43').Value = '''2.78'
$ws.Range('E43').Value = '  -2.15%  '

$ws.Range('D44').Value = '''0.0416'
$ws.Range('E44').Value = '  +2.28%  '

$ws.Range('B45').Value = 'ApeXProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D45').Value = '''3.21'
$ws.Range('E45').Value = '  +0.20%  '

$ws.Range('B46').Value = 'Stellar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D46').Value = '''0.133'
$ws.Range('E46').Value = '  +3.61%  '

$ws.Range('D47').Value = '''2.43'
$ws.Range('E47').Value = '  +0.03%  '

$ws.Range('D48').Value = '''0.999'

$ws.Range('D49').Value = '''140.64'
$ws.Range('E49').Value = '  +5.51%  '

$ws.Range('E50').Value = '  -2.09%  '

$ws.Range('D51').Value = '''8.10'
$ws.Range('E51').Value = '  -1.13%  '
